$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data: Day 2 ---
$ws.Range("A3").Value = "Day 2"

# Date cell: copy the format from B2 (date-style xf) so the existing
# numFmtId=14 style is reused instead of Excel minting a new numFmt.
$ws.Range("B3").Value = 45804
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("C3").Value = "Valid Anagram"
$ws.Range("D3").Value = "Valid Palindrome"
$ws.Range("E3").Value = "Group Anagrams"
$ws.Range("F3").Value = "HashMaps, Arrays, Two Pointers, Sorting, String"
$ws.Range("G3").Value = "S"
$ws.Range("H3").Value = "YES"

# Highlight E3 with a solid yellow fill.
$ws.Range("E3").Interior.Color = 65535  # RGB(255,255,0) -> FFFFFF00

# Widen column F to fit the longer topics text.
$ws.Columns.Item(6).ColumnWidth = 37.6666666666667

# Clear the clipboard marching ants / selection, matching the saved file.
$excel.CutCopyMode = 0
$ws.Range("C4").Select() | Out-Null
